$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.5
$ws.Range("P2").Value = 0.1
$ws.Range("S2").Value = 0.2
$ws.Range("P3").Value = 1
$ws.Range("S4").Value = 1
$ws.Range("J6").Value = 0.5555555555555556
$ws.Range("Q6").Value = 0.1111111111111111
$ws.Range("S6").Value = 0.3333333333333333
$ws.Range("B7").Value = 0.08333333333333333
$ws.Range("Q7").Value = 0.1666666666666667
$ws.Range("R7").Value = 0.08333333333333333
$ws.Range("S7").Value = 0.6666666666666666
$ws.Range("B8").Value = 0.09090909090909091
$ws.Range("F8").Value = 0.04545454545454546
$ws.Range("J8").Value = 0.09090909090909091
$ws.Range("O8").Value = 0.04545454545454546
$ws.Range("R8").Value = 0.04545454545454546
$ws.Range("S8").Value = 0.6818181818181818
$ws.Range("Q9").Value = 0.3333333333333333
$ws.Range("S9").Value = 0.6666666666666666
$ws.Range("B10").Value = 0.09302325581395349
$ws.Range("D10").Value = 0.02325581395348837
$ws.Range("F10").Value = 0.06976744186046512
$ws.Range("J10").Value = 0.1162790697674419
$ws.Range("Q10").Value = 0.2093023255813954
$ws.Range("R10").Value = 0.09302325581395349
$ws.Range("S10").Value = 0.3953488372093023
$ws.Range("G11").Value = 0.1176470588235294
$ws.Range("K11").Value = 0.2352941176470588
$ws.Range("L11").Value = 0.5882352941176471
$ws.Range("S11").Value = 0.05882352941176471
$ws.Range("G12").Value = 0.8
$ws.Range("J12").Value = 0.2
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.125
$ws.Range("H15").Value = 0.25
$ws.Range("J15").Value = 0.5
$ws.Range("S15").Value = 0.125
$ws.Range("H16").Value = 0.5
$ws.Range("I16").Value = 0.1666666666666667
$ws.Range("J16").Value = 0.1666666666666667
$ws.Range("K16").Value = 0.1666666666666667
$ws.Range("F17").Value = 0.07142857142857142
$ws.Range("H17").Value = 0.2142857142857143
$ws.Range("I17").Value = 0.1428571428571428
$ws.Range("J17").Value = 0.2857142857142857
$ws.Range("O17").Value = 0.07142857142857142
$ws.Range("S17").Value = 0.2142857142857143
$ws.Range("H18").Value = 0.1666666666666667
$ws.Range("I18").Value = 0.1666666666666667
$ws.Range("J18").Value = 0.1666666666666667
$ws.Range("K18").Value = 0.3333333333333333
$ws.Range("S18").Value = 0.1666666666666667
$ws.Range("F19").Value = 0.03225806451612903
$ws.Range("H19").Value = 0.2258064516129032
$ws.Range("I19").Value = 0.03225806451612903
$ws.Range("J19").Value = 0.2903225806451613
$ws.Range("K19").Value = 0.1612903225806452
$ws.Range("M19").Value = 0.04838709677419355
$ws.Range("O19").Value = 0.08064516129032258
$ws.Range("S19").Value = 0.1290322580645161
